$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (old row 22, timestamp 2000); this also updates the used-range dimension to A1:H21
$ws.Rows.Item(22).Delete()

# Replace sensor readings (columns C:H) for rows 2-21 with the updated dataset values
# (rows 4-21 are the previous rows 2-19 shifted down by two rows; rows 2-3 are newly added samples)
$ws.Cells.Item(2, 3).Value = -2.025566756725311
$ws.Cells.Item(2, 4).Value = 3.52062651515007
$ws.Cells.Item(2, 5).Value = 2.27691987156868
$ws.Cells.Item(2, 6).Value = -0.0397062413394451
$ws.Cells.Item(2, 7).Value = -0.0024434609804302
$ws.Cells.Item(2, 8).Value = 0.0332921557128429

$ws.Cells.Item(3, 3).Value = -2.230706214904786
$ws.Cells.Item(3, 4).Value = 3.561713695526123
$ws.Cells.Item(3, 5).Value = 2.031704187393189
$ws.Cells.Item(3, 6).Value = -0.0200058370828628
$ws.Cells.Item(3, 7).Value = -0.0035124751739203
$ws.Cells.Item(3, 8).Value = 0.0421497002243995

$ws.Cells.Item(4, 3).Value = -2.475497364997865
$ws.Cells.Item(4, 4).Value = 3.523229420185089
$ws.Cells.Item(4, 5).Value = 2.537566900253297
$ws.Cells.Item(4, 6).Value = -0.0039706239476799
$ws.Cells.Item(4, 7).Value = 0.0705549344420433
$ws.Cells.Item(4, 8).Value = 0.0328340083360672

$ws.Cells.Item(5, 3).Value = -2.48767375946045
$ws.Cells.Item(5, 4).Value = 3.061535835266112
$ws.Cells.Item(5, 5).Value = 3.86157149076462
$ws.Cells.Item(5, 6).Value = -0.0146607663482427
$ws.Cells.Item(5, 7).Value = 0.2353358417749405
$ws.Cells.Item(5, 8).Value = 0.07696902006864539

$ws.Cells.Item(6, 3).Value = -2.442349374294281
$ws.Cells.Item(6, 4).Value = 2.902286112308502
$ws.Cells.Item(6, 5).Value = 4.137762367725372
$ws.Cells.Item(6, 6).Value = 0.0704022198915481
$ws.Cells.Item(6, 7).Value = 0.2420553565025329
$ws.Cells.Item(6, 8).Value = 0.08109235763549801

$ws.Cells.Item(7, 3).Value = -2.809046030044556
$ws.Cells.Item(7, 4).Value = 1.600962877273557
$ws.Cells.Item(7, 5).Value = 4.795935153961183
$ws.Cells.Item(7, 6).Value = 0.6050620079040527
$ws.Cells.Item(7, 7).Value = 0.2802344262599945
$ws.Cells.Item(7, 8).Value = 0.07590000331401819

$ws.Cells.Item(8, 3).Value = -3.426159977912904
$ws.Cells.Item(8, 4).Value = -0.2882512211799662
$ws.Cells.Item(8, 5).Value = 5.237384021282194
$ws.Cells.Item(8, 6).Value = 0.1878410577774047
$ws.Cells.Item(8, 7).Value = 0.7900014519691467
$ws.Cells.Item(8, 8).Value = 0.304669052362442

$ws.Cells.Item(9, 3).Value = -4.615099787712106
$ws.Cells.Item(9, 4).Value = 0.4410536289215286
$ws.Cells.Item(9, 5).Value = 4.911977410316469
$ws.Cells.Item(9, 6).Value = 0.3132211565971374
$ws.Cells.Item(9, 7).Value = -0.3736968040466308
$ws.Cells.Item(9, 8).Value = 0.0494800843298435

$ws.Cells.Item(10, 3).Value = -5.570275843143422
$ws.Cells.Item(10, 4).Value = 5.792195498943334
$ws.Cells.Item(10, 5).Value = 7.507518291473446
$ws.Cells.Item(10, 6).Value = 0.6346889734268188
$ws.Cells.Item(10, 7).Value = -0.973566472530365
$ws.Cells.Item(10, 8).Value = 0.2535090744495392

$ws.Cells.Item(11, 3).Value = 4.738878250122034
$ws.Cells.Item(11, 4).Value = 7.207733154296871
$ws.Cells.Item(11, 5).Value = 21.9623451232909
$ws.Cells.Item(11, 6).Value = -2.739883422851562
$ws.Cells.Item(11, 7).Value = 3.34677791595459
$ws.Cells.Item(11, 8).Value = -2.103209018707275

$ws.Cells.Item(12, 3).Value = -3.471131086349493
$ws.Cells.Item(12, 4).Value = 6.420480489730835
$ws.Cells.Item(12, 5).Value = -4.152171969413772
$ws.Cells.Item(12, 6).Value = -3.001944541931152
$ws.Cells.Item(12, 7).Value = 0.3590360581874847
$ws.Cells.Item(12, 8).Value = 1.426217675209045

$ws.Cells.Item(13, 3).Value = -2.357546925544736
$ws.Cells.Item(13, 4).Value = 3.453876137733455
$ws.Cells.Item(13, 5).Value = 3.317040443420424
$ws.Cells.Item(13, 6).Value = -0.064446285367012
$ws.Cells.Item(13, 7).Value = -0.0032070425804704
$ws.Cells.Item(13, 8).Value = -1.773189067840576

$ws.Cells.Item(14, 3).Value = -0.1309916377067535
$ws.Cells.Item(14, 4).Value = 3.795689940452578
$ws.Cells.Item(14, 5).Value = 3.876870155334465
$ws.Cells.Item(14, 6).Value = -0.5893322229385376
$ws.Cells.Item(14, 7).Value = -0.0039706239476799
$ws.Cells.Item(14, 8).Value = -0.3778201639652252

$ws.Cells.Item(15, 3).Value = -0.2815589904785198
$ws.Cells.Item(15, 4).Value = 4.85230040550232
$ws.Cells.Item(15, 5).Value = 2.205311059951784
$ws.Cells.Item(15, 6).Value = 0.3320052623748779
$ws.Cells.Item(15, 7).Value = -0.7985535860061646
$ws.Cells.Item(15, 8).Value = 0.4050036668777466

$ws.Cells.Item(16, 3).Value = -1.535586237907411
$ws.Cells.Item(16, 4).Value = 4.95532476902008
$ws.Cells.Item(16, 5).Value = 1.937351673841474
$ws.Cells.Item(16, 6).Value = 0.0734565481543541
$ws.Cells.Item(16, 7).Value = -1.579086661338806
$ws.Cells.Item(16, 8).Value = -1.274570345878601

$ws.Cells.Item(17, 3).Value = -1.377771139144897
$ws.Cells.Item(17, 4).Value = 4.249351501464844
$ws.Cells.Item(17, 5).Value = 1.058028712868691
$ws.Cells.Item(17, 6).Value = -0.2518292069435119
$ws.Cells.Item(17, 7).Value = -0.9622654914855956
$ws.Cells.Item(17, 8).Value = -0.384845107793808

$ws.Cells.Item(18, 3).Value = 0.1868795156478872
$ws.Cells.Item(18, 4).Value = 3.352623224258432
$ws.Cells.Item(18, 5).Value = 1.056536458432678
$ws.Cells.Item(18, 6).Value = -0.0522289797663688
$ws.Cells.Item(18, 7).Value = -0.2220495194196701
$ws.Cells.Item(18, 8).Value = -0.2014328092336654

$ws.Cells.Item(19, 3).Value = -0.04321670532226736
$ws.Cells.Item(19, 4).Value = 5.682518005371086
$ws.Cells.Item(19, 5).Value = 1.995282649993894
$ws.Cells.Item(19, 6).Value = 0.0580321997404098
$ws.Cells.Item(19, 7).Value = -0.2347249686717987
$ws.Cells.Item(19, 8).Value = 0.4702135324478149

$ws.Cells.Item(20, 3).Value = -0.4393689632415791
$ws.Cells.Item(20, 4).Value = 3.92675977945327
$ws.Cells.Item(20, 5).Value = 1.474137753248212
$ws.Cells.Item(20, 6).Value = 0.08491026610136029
$ws.Cells.Item(20, 7).Value = 0.1505782902240753
$ws.Cells.Item(20, 8).Value = -0.0226020142436027

$ws.Cells.Item(21, 3).Value = -0.5461759567260746
$ws.Cells.Item(21, 4).Value = 4.010827064514163
$ws.Cells.Item(21, 5).Value = 1.266485691070556
$ws.Cells.Item(21, 6).Value = 0.0704022198915481
$ws.Cells.Item(21, 7).Value = 0.0319177098572254
$ws.Cells.Item(21, 8).Value = 0.1357648074626922

